# Automatische test-sync: 2025-07-22 12:30:50
#
# Adds a new mail-log entry (row 9) to the "Logs" sheet and updates the
# "Dashboard" sheet summary table so that the re-sorted category counts
# are reflected (Productinformatie now ties with Retour / Terugbetaling
# at 3 occurrences each).

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 9 -------------------------------------------------

$logs.Range("A9").Value = "Ik wil een handleiding ontvangen voor model EcoPro-700."
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Testmail #9: Ik wil een handleiding ontvangen voor model EcoPro-700."
$logs.Range("D9").Value = "Productinformatie"
$logs.Range("E9").Value = "Beste klant,`nBedankt voor uw bericht. Helaas hebben wij geen informatie over een model genaamd EcoPro-700 in onze systemen. Kunt u ons wat meer details geven over het product of de fabrikant, zodat we u verder kunnen helpen?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F9").Value = "2025-07-22 12:29:58"
$logs.Range("G9").Value = "Ja"
$logs.Range("H9").Value = "Nee"
$logs.Range("I9").Value = "Ja"
$logs.Range("J9").Value = "Ja"

# --- Logs sheet: extend conditional formatting ranges to include row 9 -------

foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`8")
    $newRange = $logs.Range("$col`2:$col`9")
    foreach ($fc in $oldRange.FormatConditions) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: re-sort categories now that Productinformatie ties -----

$dashboard.Range("A2").Value = "Productinformatie"
$dashboard.Range("B2").Value = 3
$dashboard.Range("A3").Value = "Retour / Terugbetaling"
$dashboard.Range("B3").Value = 3
